# Apply updated market/profit values per sheet, as captured by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 801544.1
$ws.Range("I6").Value = 1488185
$ws.Range("J6").Value = 463.16666
$ws.Range("K6").Value = 4464555
$ws.Range("L6").Value = 1389.49998
$ws.Range("M6").Value = -4464443
$ws.Range("N6").Value = -1613.49998
$ws.Range("H70").Value = 2484.375
$ws.Range("I70").Value = 3121
$ws.Range("J70").Value = 1423.3334
$ws.Range("K70").Value = 9363
$ws.Range("L70").Value = 4270.0002
$ws.Range("M70").Value = -9093
$ws.Range("N70").Value = -4810.0002
$ws.Range("H73").Value = 2484.375
$ws.Range("I73").Value = 3121
$ws.Range("J73").Value = 1423.3334
$ws.Range("K73").Value = 9363
$ws.Range("L73").Value = 4270.0002
$ws.Range("M73").Value = -8427
$ws.Range("N73").Value = -6142.0002
$ws.Range("H107").Value = 802
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 903
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 903
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -4743
$ws.Range("H112").Value = 19613.023
$ws.Range("I112").Value = 369.25
$ws.Range("J112").Value = 21586.744
$ws.Range("K112").Value = 1107.75
$ws.Range("L112").Value = 64760.232
$ws.Range("M112").Value = 0.25
$ws.Range("N112").Value = -66976.23199999999
$ws.Range("H129").Value = 1033.4807
$ws.Range("J129").Value = 1011.2195
$ws.Range("L129").Value = 3033.6585
$ws.Range("N129").Value = -13033.6585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14494961
$ws.Range("I132").Value = 1964
$ws.Range("K132").Value = 5892
$ws.Range("M132").Value = -3362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.833336
$ws.Range("I2").Value = 14.545455
$ws.Range("J2").Value = 146.42857
$ws.Range("K2").Value = 87.27273
$ws.Range("L2").Value = 878.57142
$ws.Range("M2").Value = 25.72727
$ws.Range("N2").Value = -1104.57142
$ws.Range("H107").Value = 48489740
$ws.Range("I107").Value = 233.75
$ws.Range("J107").Value = 76198030
$ws.Range("K107").Value = 701.25
$ws.Range("L107").Value = 228594090
$ws.Range("M107").Value = 1218.75
$ws.Range("N107").Value = -228597930
$ws.Range("H122").Value = 48085530
$ws.Range("I122").Value = 125000200
$ws.Range("J122").Value = 13856.125
$ws.Range("K122").Value = 1125001800
$ws.Range("L122").Value = 124705.125
$ws.Range("M122").Value = -1124999350
$ws.Range("N122").Value = -129605.125
$ws.Range("H131").Value = 809.26
$ws.Range("J131").Value = 821.0928
$ws.Range("L131").Value = 2463.2784
$ws.Range("N131").Value = -12543.2784
$ws.Range("H132").Value = 13892908
$ws.Range("I132").Value = 650
$ws.Range("J132").Value = 20005502
$ws.Range("K132").Value = 5850
$ws.Range("L132").Value = 180049518
$ws.Range("M132").Value = -3320
$ws.Range("N132").Value = -180054578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16671806
$ws.Range("I80").Value = 6166.6
$ws.Range("J80").Value = 100000000
$ws.Range("K80").Value = 6166.6
$ws.Range("L80").Value = 100000000
$ws.Range("M80").Value = -5168.6
$ws.Range("N80").Value = -100001996
$ws.Range("H83").Value = 16671806
$ws.Range("I83").Value = 6166.6
$ws.Range("J83").Value = 100000000
$ws.Range("K83").Value = 30833
$ws.Range("L83").Value = 500000000
$ws.Range("M83").Value = -25841
$ws.Range("N83").Value = -500009984
$ws.Range("H122").Value = 20006100
$ws.Range("I122").Value = 27785004
$ws.Range("J122").Value = 3202.2856
$ws.Range("K122").Value = 83355012
$ws.Range("L122").Value = 9606.856800000001
$ws.Range("M122").Value = -83352562
$ws.Range("N122").Value = -14506.8568
$ws.Range("H132").Value = 4370.6387
$ws.Range("I132").Value = 1420.125
$ws.Range("J132").Value = 27974.75
$ws.Range("K132").Value = 4260.375
$ws.Range("L132").Value = 83924.25
$ws.Range("M132").Value = -1730.375
$ws.Range("N132").Value = -88984.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1611.125
$ws.Range("I7").Value = 1478
$ws.Range("J7").Value = 1833
$ws.Range("K7").Value = 1478
$ws.Range("L7").Value = 1833
$ws.Range("M7").Value = -1366
$ws.Range("N7").Value = -2057
$ws.Range("H40").Value = 35719870
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 50007100
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 50007100
$ws.Range("M40").Value = -1664
$ws.Range("N40").Value = -50007372
$ws.Range("H122").Value = 8078.95
$ws.Range("I122").Value = 10931.667
$ws.Range("J122").Value = 3799.875
$ws.Range("K122").Value = 32795.001
$ws.Range("L122").Value = 11399.625
$ws.Range("M122").Value = -30345.001
$ws.Range("N122").Value = -16299.625
$ws.Range("H126").Value = 1611.125
$ws.Range("I126").Value = 1478
$ws.Range("J126").Value = 1833
$ws.Range("K126").Value = 4434
$ws.Range("L126").Value = 5499
$ws.Range("M126").Value = -1964
$ws.Range("N126").Value = -10439

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4471.4287
$ws.Range("I62").Value = 3660
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 3660
$ws.Range("L62").Value = 6500
$ws.Range("M62").Value = -3036
$ws.Range("N62").Value = -7748
$ws.Range("H64").Value = 17754.8
$ws.Range("J64").Value = 17754.8
$ws.Range("L64").Value = 17754.8
$ws.Range("N64").Value = -18250.8
$ws.Range("H65").Value = 4471.4287
$ws.Range("I65").Value = 3660
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 18300
$ws.Range("L65").Value = 32500
$ws.Range("M65").Value = -15180
$ws.Range("N65").Value = -38740
$ws.Range("H67").Value = 17754.8
$ws.Range("J67").Value = 17754.8
$ws.Range("L67").Value = 17754.8
$ws.Range("N67").Value = -19470.8
$ws.Range("H126").Value = 1446.7693
$ws.Range("I126").Value = 1000.8
$ws.Range("J126").Value = 2933.3333
$ws.Range("K126").Value = 3002.4
$ws.Range("L126").Value = 8799.999899999999
$ws.Range("M126").Value = -532.3999999999996
$ws.Range("N126").Value = -13739.9999
$ws.Range("H132").Value = 20904.947
$ws.Range("I132").Value = 26029.05
$ws.Range("J132").Value = 8094.6875
$ws.Range("K132").Value = 78087.14999999999
$ws.Range("L132").Value = 24284.0625
$ws.Range("M132").Value = -75557.14999999999
$ws.Range("N132").Value = -29344.0625
